# Corrected some selection scopes
# Removes rows that were incorrectly included in the naive-forecaster
# qoq error series (every other quarter in the first half of the
# history was duplicated/misaligned), shifting the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to remove (1-based, referring to the ORIGINAL row numbers before
# any deletions). Delete from the bottom up so earlier row numbers stay
# valid as we go.
$rowsToRemove = @(30, 28, 26, 24, 22, 20, 18, 16, 14, 12, 10, 8, 6, 4, 2)

foreach ($r in $rowsToRemove) {
    $ws.Rows.Item($r).Delete()
}
